$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "2026-01-28 04:01"
$ws.Range("B3").Value = 39
$ws.Range("C3").Value = 7
